$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 143
$ws.Range("F3").Value = 344
$ws.Range("F4").Value = 448
$ws.Range("F5").Value = 1761
$ws.Range("F6").Value = 90
$ws.Range("F7").Value = 2222
$ws.Range("F8").Value = 9
$ws.Range("F9").Value = 285
$ws.Range("F10").Value = 107
$ws.Range("F11").Value = 5018
$ws.Range("F12").Value = 66
$ws.Range("F13").Value = 38
$ws.Range("F14").Value = 311
$ws.Range("F15").Value = 233
$ws.Range("F16").Value = 34
$ws.Range("F17").Value = 192
$ws.Range("F18").Value = 250
$ws.Range("F19").Value = 22
$ws.Range("F20").Value = 126
$ws.Range("F21").Value = 4010
$ws.Range("F22").Value = 722
$ws.Range("F23").Value = 704
$ws.Range("F24").Value = 31
$ws.Range("F25").Value = 25
$ws.Range("F26").Value = 113
$ws.Range("F27").Value = 127
$ws.Range("F29").Value = 13
$ws.Range("F31").Value = 589
$ws.Range("F32").Value = 12
$ws.Range("F33").Value = 28
$ws.Range("F34").Value = 1025
$ws.Range("F35").Value = 5
$ws.Range("F36").Value = 2605
$ws.Range("F37").Value = 431
$ws.Range("F38").Value = 32

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 1

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 143
$ws.Range("F3").Value = 344
$ws.Range("F4").Value = 448
$ws.Range("F5").Value = 1761
$ws.Range("F6").Value = 90
$ws.Range("F7").Value = 2222
$ws.Range("F8").Value = 9
$ws.Range("F9").Value = 285
$ws.Range("F10").Value = 107
$ws.Range("F11").Value = 5018
$ws.Range("F12").Value = 66
$ws.Range("F13").Value = 38
$ws.Range("F14").Value = 311
$ws.Range("F15").Value = 233
$ws.Range("F16").Value = 34
$ws.Range("F17").Value = 192
$ws.Range("F18").Value = 252
$ws.Range("F19").Value = 22
$ws.Range("F20").Value = 126
$ws.Range("F21").Value = 4010
$ws.Range("F22").Value = 722
$ws.Range("F23").Value = 704
$ws.Range("F24").Value = 31
$ws.Range("F25").Value = 25
$ws.Range("F27").Value = 127
$ws.Range("F28").Value = 27
$ws.Range("F29").Value = 13
$ws.Range("F30").Value = 95
$ws.Range("F31").Value = 589
$ws.Range("F32").Value = 12
$ws.Range("F33").Value = 1
$ws.Range("F35").Value = 1025
$ws.Range("F36").Value = 5
$ws.Range("F37").Value = 2605
$ws.Range("F38").Value = 431
$ws.Range("F39").Value = 32
